# Update cryptocurrency Price (D) and Volume(1h) (E) figures for Sheet1,
# matching the GitHub Actions data refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.976.58'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.561.47'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.45'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.12'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').Value = '1.784.15'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = '1.563.18'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '61.89'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '26.958.82'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.68'
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.93'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.35'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.10'
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.12'
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('D34').Value = '1.422.60'
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('E35').Value = '  +2.30%  '
$ws.Range('E36').Value = '  +8.33%  '
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('E39').Value = '  +2.43%  '
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.54'
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').Value = '1.697.45'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.24'
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0957'
$ws.Range('E51').Value = '  +0.29%  '
